$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "25.302.28", "  -2.33%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.562.74", "  -3.56%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.40%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "207.24", "  -2.94%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  -0.39%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.478", "  -4.33%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.0613", "  -0.52%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.244", "  -2.25%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "17.94", "  -1.66%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.0782", "  -1.01%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.778.55", "  -3.66%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.01", "  -3.66%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.549.67", "  -4.46%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.509", "  -2.99%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "25.310.92", "  -2.32%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0₃0715", "  -3.04%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "59.31", "  -3.55%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  -0.43%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "187.08", "  -2.44%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.14", "  -2.24%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "9.34", "  -1.65%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "5.88", "  -2.40%  "),
    @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.130", "  -2.64%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "140.57", "  -2.26%  "),
    @("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.00", "  -0.41%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.67", "  -3.15%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "14.97", "  -1.96%  "),
    @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "6.43", "  -3.69%  "),
    @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.15", "  -6.56%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0468", "  -2.81%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.06", "  -2.15%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "2.99", "  -3.66%  "),
    @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.48", "  -0.63%  "),
    @("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.32", "  -3.93%  "),
    @("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.084.79", "  -3.23%  "),
    @("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.00", "  -0.64%  "),
    @("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.31", "  -3.61%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0149", "  -2.18%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.497", "  -3.66%  "),
    @("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.773", "  -8.73%  "),
    @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.804", "  +7.28%  "),
    @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "93.10", "  -4.91%  "),
    @("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "5.09", "  -0.19%  "),
    @("RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "1.693.71", "  -3.57%  "),
    @("BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.0₆0112", "  -1.10%  "),
    @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.47", "  -1.13%  "),
    @("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "52.65", "  -2.72%  "),
    @("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.0505", "  -3.40%  "),
    @("Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.405", "  -1.87%  "),
    @("USDD", "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd", "1.00", "  -0.55%  ")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $cB = $ws.Cells.Item($r, 2)
    $cC = $ws.Cells.Item($r, 3)
    $cD = $ws.Cells.Item($r, 4)
    $cE = $ws.Cells.Item($r, 5)

    # Columns D and E hold text-formatted values (e.g. "1.00", "  -2.33%  ")
    # that must not be auto-converted to numbers by Excel.
    $cD.NumberFormat = "@"
    $cE.NumberFormat = "@"

    $cB.Value = $row[0]
    $cC.Value = $row[1]
    $cD.Value = $row[2]
    $cE.Value = $row[3]
}
